# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on a handful of rows across several of the
# per-job "Leve" profit sheets, reflecting the latest Market Board pricing.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 370.44
$ws.Range("I53").Value = 400.7
$ws.Range("J53").Value = 350.26666
$ws.Range("K53").Value = 400.7
$ws.Range("L53").Value = 350.26666
$ws.Range("M53").Value = 236.3
$ws.Range("N53").Value = -1624.26666

$ws.Range("H64").Value = 3065
$ws.Range("I64").Value = 3106.0715
$ws.Range("J64").Value = 2982.8572
$ws.Range("K64").Value = 3106.0715
$ws.Range("L64").Value = 2982.8572
$ws.Range("M64").Value = -2858.0715
$ws.Range("N64").Value = -3478.8572

$ws.Range("H67").Value = 3065
$ws.Range("I67").Value = 3106.0715
$ws.Range("J67").Value = 2982.8572
$ws.Range("K67").Value = 3106.0715
$ws.Range("L67").Value = 2982.8572
$ws.Range("M67").Value = -2248.0715
$ws.Range("N67").Value = -4698.8572

$ws.Range("H113").Value = 3365.6875
$ws.Range("I113").Value = 3480.5
$ws.Range("J113").Value = 3174.3333
$ws.Range("K113").Value = 3480.5
$ws.Range("L113").Value = 3174.3333
$ws.Range("M113").Value = -226.5
$ws.Range("N113").Value = -9682.3333

$ws.Range("H127").Value = 584.4286
$ws.Range("I127").Value = 535.2
$ws.Range("K127").Value = 1605.6
$ws.Range("M127").Value = 3354.4

$ws.Range("H141").Value = 423674.25
$ws.Range("I141").Value = 2375.111
$ws.Range("J141").Value = 564107.3
$ws.Range("K141").Value = 7125.333
$ws.Range("L141").Value = 1692321.9
$ws.Range("M141").Value = -1945.333
$ws.Range("N141").Value = -1702681.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5826.35
$ws.Range("I32").Value = 5038.3516
$ws.Range("J32").Value = 13793.889
$ws.Range("K32").Value = 5038.3516
$ws.Range("L32").Value = 13793.889
$ws.Range("M32").Value = -4751.3516
$ws.Range("N32").Value = -14367.889

$ws.Range("H97").Value = 915.16
$ws.Range("I97").Value = 874.95
$ws.Range("J97").Value = 1076
$ws.Range("K97").Value = 874.95
$ws.Range("L97").Value = 1076
$ws.Range("M97").Value = -378.95
$ws.Range("N97").Value = -2068

$ws.Range("H110").Value = 1432.0416
$ws.Range("I110").Value = 538.45
$ws.Range("K110").Value = 538.45
$ws.Range("M110").Value = 1506.55

$ws.Range("H122").Value = 2201.303
$ws.Range("I122").Value = 1676.1154
$ws.Range("J122").Value = 4152
$ws.Range("K122").Value = 5028.3462
$ws.Range("L122").Value = 12456
$ws.Range("M122").Value = -2578.3462
$ws.Range("N122").Value = -17356

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2102.8333
$ws.Range("I20").Value = 1699.5
$ws.Range("K20").Value = 1699.5
$ws.Range("M20").Value = -1452.5

$ws.Range("H94").Value = 719.75
$ws.Range("I94").Value = 719.75
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 719.75
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -268.75
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 24933.572
$ws.Range("I4").Value = 4882.75
$ws.Range("J4").Value = 51668
$ws.Range("K4").Value = 4882.75
$ws.Range("L4").Value = 51668
$ws.Range("M4").Value = -4770.75
$ws.Range("N4").Value = -51892

$ws.Range("H5").Value = 58097700
$ws.Range("I5").Value = 141093710
$ws.Range("J5").Value = 494.8
$ws.Range("K5").Value = 141093710
$ws.Range("L5").Value = 494.8
$ws.Range("M5").Value = -141093598
$ws.Range("N5").Value = -718.8

$ws.Range("H8").Value = 803
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 954.5
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 954.5
$ws.Range("M8").Value = -360
$ws.Range("N8").Value = -1234.5

$ws.Range("H11").Value = 2900
$ws.Range("I11").Value = 2900
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2900
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -2760
$ws.Range("N11").ClearContents()

$ws.Range("H31").Value = 2502837.8
$ws.Range("I31").Value = 2780020.2
$ws.Range("J31").Value = 8195
$ws.Range("K31").Value = 2780020.2
$ws.Range("L31").Value = 8195
$ws.Range("M31").Value = -2779725.2
$ws.Range("N31").Value = -8785

$ws.Range("H34").Value = 2502837.8
$ws.Range("I34").Value = 2780020.2
$ws.Range("J34").Value = 8195
$ws.Range("K34").Value = 2780020.2
$ws.Range("L34").Value = 8195
$ws.Range("M34").Value = -2779818.2
$ws.Range("N34").Value = -8599

$ws.Range("H53").Value = 23817
$ws.Range("J53").Value = 23817
$ws.Range("L53").Value = 23817
$ws.Range("N53").Value = -25031

$ws.Range("H99").Value = 1851.0555
$ws.Range("I99").Value = 1032.6
$ws.Range("K99").Value = 1032.6
$ws.Range("M99").Value = 465.4000000000001

$ws.Range("H126").Value = 1851.0555
$ws.Range("I126").Value = 1032.6
$ws.Range("K126").Value = 3097.8
$ws.Range("M126").Value = -627.7999999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8949.454
$ws.Range("I3").Value = 680
$ws.Range("J3").Value = 12050.5
$ws.Range("K3").Value = 680
$ws.Range("L3").Value = 12050.5
$ws.Range("M3").Value = -564
$ws.Range("N3").Value = -12282.5

$ws.Range("H4").Value = 83337.336
$ws.Range("J4").Value = 83337.336
$ws.Range("L4").Value = 83337.336
$ws.Range("N4").Value = -83561.336

$ws.Range("H70").Value = 7336
$ws.Range("I70").Value = 8008
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 8008
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -7738
$ws.Range("N70").Value = -7540

$ws.Range("H73").Value = 7336
$ws.Range("I73").Value = 8008
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 8008
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -7072
$ws.Range("N73").Value = -8872

$ws.Range("H102").Value = 33111.395
$ws.Range("I102").Value = 2085.7727
$ws.Range("J102").Value = 95162.63
$ws.Range("K102").Value = 2085.7727
$ws.Range("L102").Value = 95162.63
$ws.Range("M102").Value = -463.7727
$ws.Range("N102").Value = -98406.63

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5264905.5
$ws.Range("I7").Value = 10000970
$ws.Range("J7").Value = 2611.6667
$ws.Range("K7").Value = 10000970
$ws.Range("L7").Value = 2611.6667
$ws.Range("M7").Value = -10000858
$ws.Range("N7").Value = -2835.6667

$ws.Range("H40").Value = 3800.4707
$ws.Range("I40").Value = 3959
$ws.Range("J40").Value = 3420
$ws.Range("K40").Value = 3959
$ws.Range("L40").Value = 3420
$ws.Range("M40").Value = -3823
$ws.Range("N40").Value = -3692

$ws.Range("H46").Value = 1313.5416
$ws.Range("I46").Value = 954.87805
$ws.Range("K46").Value = 954.87805
$ws.Range("M46").Value = -766.87805

$ws.Range("H93").Value = 3251
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 4001.3333
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 4001.3333
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -6497.3333

$ws.Range("H126").Value = 5264905.5
$ws.Range("I126").Value = 10000970
$ws.Range("J126").Value = 2611.6667
$ws.Range("K126").Value = 30002910
$ws.Range("L126").Value = 7835.000100000001
$ws.Range("M126").Value = -30000440
$ws.Range("N126").Value = -12775.0001

$ws.Range("H135").Value = 39166.668
$ws.Range("J135").Value = 39166.668
$ws.Range("L135").Value = 39166.668
$ws.Range("N135").Value = -49306.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 251816.38
$ws.Range("I122").Value = 436163.7
$ws.Range("J122").Value = 2405.2942
$ws.Range("K122").Value = 1308491.1
$ws.Range("L122").Value = 7215.882599999999
$ws.Range("M122").Value = -1306041.1
$ws.Range("N122").Value = -12115.8826

$ws.Range("H126").Value = 3573535
$ws.Range("I126").Value = 1340.4736
$ws.Range("J126").Value = 11114834
$ws.Range("K126").Value = 4021.4208
$ws.Range("L126").Value = 33344502
$ws.Range("M126").Value = -1551.4208
$ws.Range("N126").Value = -33349442
